# Fix an error in the Bant Spirit features data (Features_old.xlsx, "Feuil1"):
# noble_hierarch (row 18) actually produces W/U/G mana, while
# mausoleum_wanderer (row 19) does not - the produces_W/produces_U/produces_G
# columns (H:J) had these swapped.
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Feuil1")
$ws.Activate()

# noble_hierarch: set produces_W, produces_U, produces_G to 1
$ws.Range("H18:J18").Value = 1

# mausoleum_wanderer: set produces_W, produces_U, produces_G to 0
$ws.Range("H19:J19").Value = 0

$ws.Range("J18").Select()
